$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text while we assign numeric-looking strings,
# matching the original inline-string cell type. ClearFormats afterwards removes
# the temporary style index so cells end up with no explicit "s" attribute, same
# as the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "23.002.31"
$ws.Range("E2").Value = "  -4.06%  "

# Row 3
$ws.Range("D3").Value = "1.599.03"
$ws.Range("E3").Value = "  -3.29%  "

# Row 4
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.55%  "

# Row 5
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").Value = "300.67"
$ws.Range("E6").Value = "  -2.63%  "

# Row 7
$ws.Range("D7").Value = "0.3758"
$ws.Range("E7").Value = "  -3.76%  "

# Row 8
$ws.Range("D8").Value = "0.3649"
$ws.Range("E8").Value = "  -4.75%  "

# Row 9
$ws.Range("D9").Value = "48.27"
$ws.Range("E9").Value = "  -6.03%  "

# Row 10
$ws.Range("D10").Value = "1.005"
$ws.Range("E10").Value = "  +0.58%  "

# Row 11
$ws.Range("D11").Value = "1.271"
$ws.Range("E11").Value = "  -6.21%  "

# Row 12
$ws.Range("D12").Value = "0.08045"
$ws.Range("E12").Value = "  -4.74%  "

# Row 13
$ws.Range("D13").Value = "22.87"
$ws.Range("E13").Value = "  -4.70%  "

# Row 14
$ws.Range("D14").Value = "6.597"
$ws.Range("E14").Value = "  -7.39%  "

# Row 15
$ws.Range("D15").Value = "7.642"
$ws.Range("E15").Value = "  -2.90%  "

# Row 16
$ws.Range("D16").Value = "0.00001260"
$ws.Range("E16").Value = "  -4.10%  "

# Row 17
$ws.Range("D17").Value = "1.592.16"
$ws.Range("E17").Value = "  -3.36%  "

# Row 18
$ws.Range("D18").Value = "91.29"
$ws.Range("E18").Value = "  -3.37%  "

# Row 19
$ws.Range("D19").Value = "0.06783"
$ws.Range("E19").Value = "  -2.67%  "

# Row 20
$ws.Range("D20").Value = "18.30"
$ws.Range("E20").Value = "  -7.45%  "

# Row 21
$ws.Range("D21").Value = "6.566"
$ws.Range("E21").Value = "  -4.91%  "

# Row 22
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("D23").Value = "12.95"
$ws.Range("E23").Value = "  -4.96%  "

# Row 24
$ws.Range("D24").Value = "23.033.96"
$ws.Range("E24").Value = "  -3.88%  "

# Row 25
$ws.Range("D25").Value = "2.355"
$ws.Range("E25").Value = "  -5.06%  "

# Row 26
$ws.Range("D26").Value = "2.907"
$ws.Range("E26").Value = "  -3.81%  "

# Row 27
$ws.Range("D27").Value = "21.08"
$ws.Range("E27").Value = "  -4.42%  "

# Row 28
$ws.Range("D28").Value = "150.30"
$ws.Range("E28").Value = "  -1.80%  "

# Row 29
$ws.Range("D29").Value = "5.267"
$ws.Range("E29").Value = "  -3.22%  "

# Row 30
$ws.Range("D30").Value = "132.15"
$ws.Range("E30").Value = "  -5.06%  "

# Row 31
$ws.Range("D31").Value = "2.422"
$ws.Range("E31").Value = "  -2.40%  "

# Row 32
$ws.Range("D32").Value = "6.914"
$ws.Range("E32").Value = "  -10.68%  "

# Row 33
$ws.Range("D33").Value = "1.772.73"
$ws.Range("E33").Value = "  -3.34%  "

# Row 34
$ws.Range("D34").Value = "0.9851"
$ws.Range("E34").Value = "  -5.00%  "

# Row 35
$ws.Range("D35").Value = "0.07686"
$ws.Range("E35").Value = "  -5.38%  "

# Row 36
$ws.Range("D36").Value = "0.02771"
$ws.Range("E36").Value = "  -6.72%  "

# Row 37
$ws.Range("D37").Value = "6.251"
$ws.Range("E37").Value = "  -7.51%  "

# Row 38
$ws.Range("D38").Value = "0.2530"
$ws.Range("E38").Value = "  -5.51%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.08845"
$ws.Range("E39").Value = "  -3.33%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "10.04"
$ws.Range("E40").Value = "  -7.44%  "

# Row 41
$ws.Range("D41").Value = "1.391"
$ws.Range("E41").Value = "  -2.49%  "

# Row 42
$ws.Range("D42").Value = "0.7118"
$ws.Range("E42").Value = "  -5.87%  "

# Row 43
$ws.Range("D43").Value = "12.71"
$ws.Range("E43").Value = "  -5.86%  "

# Row 44
$ws.Range("D44").Value = "15.86"
$ws.Range("E44").Value = "  -3.12%  "

# Row 45
$ws.Range("D45").Value = "0.6579"
$ws.Range("E45").Value = "  -5.30%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.298"
$ws.Range("E47").Value = "  -6.43%  "

# Row 48
$ws.Range("D48").Value = "3.958"
$ws.Range("E48").Value = "  -3.05%  "

# Row 49
$ws.Range("D49").Value = "0.07987"
$ws.Range("E49").Value = "  -3.78%  "

# Row 50
$ws.Range("D50").Value = "131.03"
$ws.Range("E50").Value = "  -2.59%  "

# Row 51
$ws.Range("D51").Value = "1.164"
$ws.Range("E51").Value = "  -4.89%  "

# Remove the temporary text-format style from column D so the XML has no
# leftover "s" attribute on these cells (mirrors source formatting).
$ws.Range("D2:D51").ClearFormats()
